$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").Value = "Finnish"
$ws.Range("A2").Value = "129° - Weetabix cereal 48 pack £3.50 at Asda"
$ws.Range("B2").Value = "Finnish"
$ws.Range("C2").Value = "129 ° - Weetabix-muroja 48 kpl 3,50 puntaa Asdassa"
$ws.Range("D2").Value = "https://www.hotukdeals.com/deals/rollback-ps350-weetabix-cereal-at-asda-3571010"
